$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark additional cells with the "x" string, matching the existing E5 cell.
$ws.Range("I1").Value = "x"
$ws.Range("B2").Value = "x"
$ws.Range("E2").Value = "x"
$ws.Range("G4").Value = "x"
$ws.Range("F6").Value = "x"

# Update the selection to I1 as shown in the diff.
$ws.Range("I1").Select()
